# Atualização planilha de orçamento
# Adds a "Valor Mensal Previsto" label (and adjacent formatted cells) to the
# budget row, extending the used range from A1:E4 to A1:F4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell: D4 holds the new header label, bold Arial 12, vertically centered
# (same look as the worksheet title in row 1).
$ws.Range("D4").Value = "Valor Mensal Previsto"
$ws.Range("D4").Font.Name = "Arial"
$ws.Range("D4").Font.Size = 12
$ws.Range("D4").Font.Bold = $true
$ws.Range("D4").HorizontalAlignment = 1
$ws.Range("D4").VerticalAlignment = -4108

# E4 carries the same formatting as D4 (bold, vertically centered).
$ws.Range("E4").Font.Name = "Arial"
$ws.Range("E4").Font.Size = 12
$ws.Range("E4").Font.Bold = $true
$ws.Range("E4").HorizontalAlignment = 1
$ws.Range("E4").VerticalAlignment = -4108

# F4 carries the same formatting as B4:C4 (not bold, centered both ways).
$ws.Range("F4").Font.Name = "Arial"
$ws.Range("F4").Font.Size = 12
$ws.Range("F4").Font.Bold = $false
$ws.Range("F4").HorizontalAlignment = -4108
$ws.Range("F4").VerticalAlignment = -4108

# Reflect the new used range / last-selected cell, matching the saved file.
$ws.Range("G13").Select() | Out-Null
